$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1612"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 781242

$ws.Range("E17").Value = "1712"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 781242

$ws.Range("E18").Value = "1801"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 781242

$ws.Range("E19").Value = "1802"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 781242

$ws.Range("E20").Value = "1803"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 781242

$ws.Range("E21").Value = "1804"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 781242
